$wb = $excel.ActiveWorkbook

# The "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md" row (row 4 on each sheet) just
# had its handoff report regenerated, so its timestamps advance while the
# "6edcb7bb-6a62-49b9-b5d2-6afc280672df.md" row (row 5) keeps its prior values.

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 4.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-12-16 08:12:14"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for row 4.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-12-16 08:12:01"

# de-de sheet: "Latest Handoff Datetime" column (H) for row 4.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-12-16 08:12:14"
